$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Remove the existing hyperlinks first. The rows are being reshuffled and two
# links (Library Management System project + My Portfolio) need to end up
# pointing at different cells/targets, so it's cleanest to drop every
# hyperlink and recreate them once all the text is in its final place.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# Write the cell text directly into its final resting place. The order below
# matters: it reproduces the order in which new strings were introduced into
# the shared-string table (Hospital Management System / its link first, then
# the Library-Management-System link before its label).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Hospital Management System"
$ws.Range("B1").Value = "https://github.com/Lohi26/miniproject/tree/main/717821p231"

$ws.Range("B2").Value = "https://github.com/Lohi26/Library-Management-System"
$ws.Range("A2").Value = "Online Library Management System"

$ws.Range("A3").Value = "Bankist Application"
$ws.Range("B3").Value = "https://bankist-one-roan.vercel.app/"

$ws.Range("A4").Value = "Digital CV"
$ws.Range("B4").Value = "https://my-resume-five-sand.vercel.app/"

$ws.Range("A5").Value = "Guess My Number Game"
$ws.Range("B5").Value = "https://guess-my-number-gamma-lake.vercel.app/"

$ws.Range("A6").Value = "Pig Game (Dice Roll)"
$ws.Range("B6").Value = "https://pig-game-sage-chi.vercel.app/"

$ws.Range("A7").Value = "My Portfolio"
$ws.Range("B7").Value = "https://github.com/Lohi26/MyPortfolio"

# ---------------------------------------------------------------------------
# Recreate the hyperlinks against their final cells / URLs.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B3"), "https://bankist-one-roan.vercel.app/")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://my-resume-five-sand.vercel.app/")
$ws.Hyperlinks.Add($ws.Range("B1"), "https://github.com/Lohi26/miniproject/tree/main/717821p231")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://guess-my-number-gamma-lake.vercel.app/")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://pig-game-sage-chi.vercel.app/")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/Lohi26/MyPortfolio")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/Lohi26/Library-Management-System")

# Make sure the whole link column keeps the built-in Hyperlink look (adding a
# link via the API above stamps its own one-off style, so reapply the shared
# named style across the column).
$ws.Range("B1:B7").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Column widths widened slightly to fit the new, longer text/links.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 48.17
$ws.Columns.Item(2).ColumnWidth = 52.5

# Final selection left on B4 (Digital CV's link) to match the saved view.
$ws.Range("B4").Select()
